# Updated hourly forecast report ("reporthouronfcst"):
# the old "Unnamed: 0" column (B) is dropped entirely (columns C:K shift left
# to B:J) and the underlying figures are refreshed with a new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "Unnamed: 0" column - this shifts Attesa..Delta_Offerto
# (old C:K) one column to the left (new B:J), matching every header label.
$ws.Range("B:B").EntireColumn.Delete()

# Refresh the data rows (new B2:J20 = Attesa, Risposte Eff., Offerte,
# Abb sup. 14, Abb inf. 14, Short Call min 10, Cleared, 10/03/2023,
# Delta_Offerto) with the latest pull of numbers.
$rows = @(
    ,@(2, 0,329,334,4,0,3,0,318.5,4.866562009419151)
    ,@(3, 0,155,155,0,0,2,0,206,-24.75728155339806)
    ,@(4, 0,10,10,0,0,0,0,4,150)
    ,@(5, 0,153,155,2,0,0,0,106,46.22641509433962)
    ,@(6, 0,24,24,0,0,0,0,38,-36.8421052631579)
    ,@(7, 0,24,25,1,0,1,0,82,-69.51219512195121)
    ,@(8, 0,159,164,5,1,6,0,23,613.0434782608695)
    ,@(9, 0,14,17,3,0,1,0,58,-70.68965517241379)
    ,@(10, 1,154,178,23,0,2,0,379,-53.03430079155673)
    ,@(11, 0,125,126,1,0,1,0,208,-39.42307692307693)
    ,@(12, 0,217,270,11,0,2,40,395.2,-31.68016194331984)
    ,@(13, 0,11,11,0,0,0,0,302,-96.35761589403974)
    ,@(14, 0,298,314,5,1,4,10,495,-36.56565656565657)
    ,@(15, 0,87,88,1,0,1,0,127,-30.70866141732284)
    ,@(16, 0,85,101,14,2,2,0,134,-24.6268656716418)
    ,@(17, 0,32,33,0,1,1,0,77,-57.14285714285714)
    ,@(18, 0,1,1,0,0,0,0,1,0)
    ,@(19, 0,7,7,0,0,0,0,9,-22.22222222222222)
    ,@(20, 0,19,19,0,0,0,0,38,-50)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -lt $row.Length; $col++) {
        $colLetter = [char](65 + $col)
        $ws.Range("$colLetter$r").Value = $row[$col]
    }
}
